$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Insert a new blank row at position 26, shifting the existing data row
# (the "YAHYA" user record, currently row 26) down to row 27.
$ws.Rows.Item(26).Insert(-4121)  # xlShiftDown

# Update the shifted record (now row 27) with the latest user data:
# uid changed, counter incremented, roomID changed.
$ws.Range("C27").Value = "591EF2D4"
$ws.Range("D27").Value = 4
$ws.Range("E27").Value = "'105"

# The leading apostrophe above (used to keep "105" stored as text instead of
# a number) resets E27's cell style; restore it to match the rest of the row.
$ws.Range("A27").Copy() | Out-Null
$ws.Range("E27").PasteSpecial(-4122) | Out-Null  # xlPasteFormats
$excel.CutCopyMode = 0
